$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 310.75
$ws.Range("I33").Value = 254.4
$ws.Range("K33").Value = 254.4
$ws.Range("M33").Value = -25.40000000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2533.2222
$ws.Range("I43").Value = 2624.75
$ws.Range("J43").Value = 2460
$ws.Range("K43").Value = 2624.75
$ws.Range("L43").Value = 2460
$ws.Range("M43").Value = -2555.75
$ws.Range("N43").Value = -2598

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 9333.333000000001
$ws.Range("I51").Value = 9333.333000000001
$ws.Range("K51").Value = 9333.333000000001
$ws.Range("M51").Value = -8849.333000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3055.5557
$ws.Range("I62").Value = 3583.1667
$ws.Range("J62").Value = 2000.3334
$ws.Range("K62").Value = 3583.1667
$ws.Range("L62").Value = 2000.3334
$ws.Range("M62").Value = -2959.1667
$ws.Range("N62").Value = -3248.3334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 3055.5557
$ws.Range("I65").Value = 3583.1667
$ws.Range("J65").Value = 2000.3334
$ws.Range("K65").Value = 17915.8335
$ws.Range("L65").Value = 10001.667
$ws.Range("M65").Value = -14795.8335
$ws.Range("N65").Value = -16241.667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 865.8333
$ws.Range("J80").Value = 852.6923
$ws.Range("L80").Value = 2558.0769
$ws.Range("N80").Value = -4554.0769

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 865.8333
$ws.Range("J83").Value = 852.6923
$ws.Range("L83").Value = 7674.2307
$ws.Range("N83").Value = -17658.2307

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 874.75
$ws.Range("J86").Value = 999.6667
$ws.Range("L86").Value = 999.6667
$ws.Range("N86").Value = -3245.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 874.75
$ws.Range("J89").Value = 999.6667
$ws.Range("L89").Value = 4998.3335
$ws.Range("N89").Value = -16230.3335

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2261.5386
$ws.Range("I112").Value = 1300
$ws.Range("J112").Value = 2550
$ws.Range("K112").Value = 3900
$ws.Range("L112").Value = 7650
$ws.Range("M112").Value = -2792
$ws.Range("N112").Value = -9866

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2302
$ws.Range("I122").Value = 1303.4286
$ws.Range("K122").Value = 3910.2858
$ws.Range("M122").Value = -1460.2858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1570.8
$ws.Range("I132").Value = 1204.7693
$ws.Range("J132").Value = 3950
$ws.Range("K132").Value = 3614.3079
$ws.Range("L132").Value = 11850
$ws.Range("M132").Value = -1084.3079
$ws.Range("N132").Value = -16910

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").ClearContents()
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1431.6666
$ws.Range("I94").Value = 1347.5
$ws.Range("K94").Value = 1347.5
$ws.Range("M94").Value = -896.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 385000
$ws.Range("J108").Value = 385000
$ws.Range("L108").Value = 385000
$ws.Range("N108").Value = -392680

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 459.36365
$ws.Range("I107").Value = 450.8889
$ws.Range("J107").Value = 497.5
$ws.Range("K107").Value = 450.8889
$ws.Range("L107").Value = 497.5
$ws.Range("M107").Value = 1469.1111
$ws.Range("N107").Value = -4337.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1457.625
$ws.Range("I134").Value = 1310.1666
$ws.Range("K134").Value = 3930.4998
$ws.Range("M134").Value = -1395.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 440.16666
$ws.Range("I7").Value = 213.66667
$ws.Range("J7").Value = 666.6667
$ws.Range("K7").Value = 641.00001
$ws.Range("L7").Value = 2000.0001
$ws.Range("M7").Value = -529.00001
$ws.Range("N7").Value = -2224.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 674.3333
$ws.Range("I107").Value = 512
$ws.Range("J107").Value = 999
$ws.Range("K107").Value = 1536
$ws.Range("L107").Value = 2997
$ws.Range("M107").Value = 384
$ws.Range("N107").Value = -6837

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 167
$ws.Range("I107").Value = 170.8
$ws.Range("K107").Value = 170.8
$ws.Range("M107").Value = 1749.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").ClearContents()
$ws.Range("N110").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").ClearContents()
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7604.8335
$ws.Range("I40").Value = 7604.8335
$ws.Range("K40").Value = 7604.8335
$ws.Range("M40").Value = -7468.8335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3466.652
$ws.Range("I46").Value = 2894.5833
$ws.Range("K46").Value = 2894.5833
$ws.Range("M46").Value = -2706.5833

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3428.5454
$ws.Range("I122").Value = 3191
$ws.Range("J122").Value = 4497.5
$ws.Range("K122").Value = 9573
$ws.Range("L122").Value = 13492.5
$ws.Range("M122").Value = -7123
$ws.Range("N122").Value = -18392.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 34999.855
$ws.Range("J26").Value = 34999.855
$ws.Range("L26").Value = 34999.855
$ws.Range("N26").Value = -35585.855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 16750.8
$ws.Range("J41").Value = 16750.8
$ws.Range("L41").Value = 16750.8
$ws.Range("N41").Value = -17530.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("N62").Value = 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("N65").Value = 0
